$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "d.ticket[...]" placeholders to "d.tickets[...]" across rows 2-3 (A:G)
for ($r = 2; $r -le 3; $r++) {
    for ($c = 1; $c -le 7; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $text = $cell.Text
        $cell.Value = $text -replace 'd\.ticket\[', 'd.tickets['
    }
}

# Update the active selection cell
$ws.Range("F19").Select()
